$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix F3 value (was -700000, now 49950)
$ws.Range("F3").Value = 49950

# Add new row 6: Clinte
$ws.Range("A6").Value = "Clinte"
$ws.Range("B6").Value = "SP"
$ws.Range("C6").Value = "BR"
$ws.Range("D6").Value = 20000
$ws.Range("E6").Value = 35000
$ws.Range("F6").Value = -15000
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "06/05/2022"
$ws.Range("G6").ClearFormats()

# Add new row 7: Wen Tech
$ws.Range("A7").Value = "Wen Tech"
$ws.Range("B7").Value = "GO"
$ws.Range("C7").Value = "BR"
$ws.Range("D7").Value = 1000000
$ws.Range("E7").Value = 900000
$ws.Range("F7").Value = 100000
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "12/12/2021"
$ws.Range("G7").ClearFormats()
